$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add P1 = 14, Q1 = 15, matching the format of the existing header cells ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Rows 2-25: swap values in columns I, K, M, O (1 <-> 2) and add new columns P, Q = 2 ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # I
    $kVal = $ws.Cells.Item($r, 11).Value2  # K
    $mVal = $ws.Cells.Item($r, 13).Value2  # M
    $oVal = $ws.Cells.Item($r, 15).Value2  # O

    $ws.Cells.Item($r, 9).Value  = 3 - $iVal
    $ws.Cells.Item($r, 11).Value = 3 - $kVal
    $ws.Cells.Item($r, 13).Value = 3 - $mVal
    $ws.Cells.Item($r, 15).Value = 3 - $oVal

    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
